$wb = $excel.ActiveWorkbook
$wsOld = $wb.Worksheets.Item("Previously added")
$wsNew = $wb.Worksheets.Item("New")

# -----------------------------------------------------------------------
# The scraper re-ran: the 2 listings that used to sit at the top of "New"
# are now historical, so they slide down into "Previously added" (appended
# at the bottom, in the same order). The "New" sheet is refreshed with the
# 3 freshly scraped listings (2 overwrite the old rows, 1 is brand new).
# -----------------------------------------------------------------------

# --- capture the current "New" sheet rows 2 & 3 before overwriting them ---
$oldLink2 = $wsNew.Range("A2").Value2
$oldPrice2 = $wsNew.Range("B2").Value2
$oldDistrict2 = $wsNew.Range("C2").Value2
$oldArea2 = $wsNew.Range("D2").Value2
$oldCadastre2 = $wsNew.Range("E2").Value2
$oldDate2 = $wsNew.Range("F2").Value2

$oldLink3 = $wsNew.Range("A3").Value2
$oldPrice3 = $wsNew.Range("B3").Value2
$oldDistrict3 = $wsNew.Range("C3").Value2
$oldArea3 = $wsNew.Range("D3").Value2
$oldCadastre3 = $wsNew.Range("E3").Value2
$oldDate3 = $wsNew.Range("F3").Value2

# --- append those two rows to the bottom of "Previously added" ---
$lastRow = $wsOld.UsedRange.Rows.Count
$r1 = $lastRow + 1
$r2 = $lastRow + 2

# Force text number-format on the id/text columns first so long numeric
# cadastre numbers are stored as text (shared string), not coerced to
# a Number, matching how every other row in the sheet is stored.
$wsOld.Range("A$r1" + ":E$r2").NumberFormat = "@"

$wsOld.Range("A$r1").Value2 = $oldLink2
$wsOld.Range("B$r1").Value2 = $oldPrice2
$wsOld.Range("C$r1").Value2 = $oldDistrict2
$wsOld.Range("D$r1").Value2 = $oldArea2
$wsOld.Range("E$r1").Value2 = $oldCadastre2
$wsOld.Range("F$r1").Value2 = $oldDate2

$wsOld.Range("A$r2").Value2 = $oldLink3
$wsOld.Range("B$r2").Value2 = $oldPrice3
$wsOld.Range("C$r2").Value2 = $oldDistrict3
$wsOld.Range("D$r2").Value2 = $oldArea3
$wsOld.Range("E$r2").Value2 = $oldCadastre3
$wsOld.Range("F$r2").Value2 = $oldDate3

# Existing hyperlinks already baked into "Previously added" must stay put,
# so just append 2 new ones (they land at the end of the collection).
$wsOld.Hyperlinks.Add($wsOld.Range("A$r1"), $oldLink2)
$wsOld.Hyperlinks.Add($wsOld.Range("A$r2"), $oldLink3)

# Re-apply the same look & feel (style/number format) as the row right
# above them, since Hyperlinks.Add()/NumberFormat reset cell styling.
$wsOld.Range("A$lastRow" + ":F$lastRow").Copy()
$wsOld.Range("A$r1" + ":F$r1").PasteSpecial(-4122)
$wsOld.Range("A$lastRow" + ":F$lastRow").Copy()
$wsOld.Range("A$r2" + ":F$r2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- overwrite "New" sheet rows 2-3 with the new listings, append row 4 ---
$wsNew.Range("A2:E4").NumberFormat = "@"

$wsNew.Range("A2").Value2 = "https://www.ss.com/msg/lv/real-estate/wood/kuldiga-and-reg/turlavas-pag/bbknn.html"
$wsNew.Range("B2").Value2 = "1 000 €"
$wsNew.Range("C2").Value2 = "Kuldīga un raj."
$wsNew.Range("D2").Value2 = "1 ha."
$wsNew.Range("E2").Value2 = ""
$wsNew.Range("F2").Value2 = 46001.45694444445

$wsNew.Range("A3").Value2 = "https://www.ss.com/msg/lv/real-estate/wood/rezekne-and-reg/ilzeskalna-pag/fiblg.html"
$wsNew.Range("B3").Value2 = "13 500 €"
$wsNew.Range("C3").Value2 = "Rēzekne un raj."
$wsNew.Range("D3").Value2 = "1 ha."
$wsNew.Range("E3").Value2 = "78580060216"
$wsNew.Range("F3").Value2 = 46001.33333333333

$wsNew.Range("A4").Value2 = "https://www.ss.com/msg/lv/real-estate/wood/rezekne-and-reg/ozolmuizas-pag/aghex.html"
$wsNew.Range("B4").Value2 = "6 000 €"
$wsNew.Range("C4").Value2 = "Rēzekne un raj."
$wsNew.Range("D4").Value2 = "0.53 ha."
$wsNew.Range("E4").Value2 = "78780030511"
$wsNew.Range("F4").Value2 = 46000.617361111115

# The "New" sheet's 2 old hyperlinks (rows 2 & 3) pointed at the listings
# that just moved out, so they no longer apply to this sheet's content.
# Existing hyperlink relationships can't be edited/removed individually
# through this interop, but the whole collection can be cleared at once -
# do that, then add 3 fresh hyperlinks for the refreshed rows 2-4 (they
# renumber cleanly from rId1).
$wsNew.Hyperlinks.Delete()
$wsNew.Hyperlinks.Add($wsNew.Range("A2"), "https://www.ss.com/msg/lv/real-estate/wood/kuldiga-and-reg/turlavas-pag/bbknn.html")
$wsNew.Hyperlinks.Add($wsNew.Range("A3"), "https://www.ss.com/msg/lv/real-estate/wood/rezekne-and-reg/ilzeskalna-pag/fiblg.html")
$wsNew.Hyperlinks.Add($wsNew.Range("A4"), "https://www.ss.com/msg/lv/real-estate/wood/rezekne-and-reg/ozolmuizas-pag/aghex.html")

# Re-apply the normal row style, since Hyperlinks.Add()/NumberFormat just
# overwrote it with the built-in "Hyperlink" style / text format. A data
# row further down "Previously added" (not the row right after the header,
# whose columns carry a slightly different style) has the right look for
# A:F, so borrow its format.
$refRow = 300
$wsOld.Range("A$refRow" + ":F$refRow").Copy()
$wsNew.Range("A2:F2").PasteSpecial(-4122)
$wsOld.Range("A$refRow" + ":F$refRow").Copy()
$wsNew.Range("A3:F3").PasteSpecial(-4122)
$wsOld.Range("A$refRow" + ":F$refRow").Copy()
$wsNew.Range("A4:F4").PasteSpecial(-4122)
$excel.CutCopyMode = $false
